$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3382
$ws.Range("F6").Value = 7902
$ws.Range("F9").Value = 4
$ws.Range("F10").Value = 1843
$ws.Range("F12").Value = 160
$ws.Range("F13").Value = 113
$ws.Range("F14").Value = 53
$ws.Range("F16").Value = 1060
$ws.Range("F19").Value = 8536
$ws.Range("F24").Value = 1043
$ws.Range("F25").Value = 1018
$ws.Range("F27").Value = 1174
$ws.Range("F28").Value = 1073
$ws.Range("F29").Value = 589
$ws.Range("F30").Value = 28
$ws.Range("F32").Value = 6
$ws.Range("F33").Value = 110
$ws.Range("F36").Value = 464
$ws.Range("F37").Value = 371
$ws.Range("F38").Value = 3513
$ws.Range("F39").Value = 929
$ws.Range("F42").Value = 504
$ws.Range("F43").Value = 102
$ws.Range("F45").Value = 639
$ws.Range("F46").Value = 46
$ws.Range("F48").Value = 22

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 340
$ws.Range("F10").Value = 200
$ws.Range("F15").Value = 2
$ws.Range("F20").Value = 40
$ws.Range("F23").Value = 109
$ws.Range("F24").Value = 6919
$ws.Range("F33").Value = 65

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2067
$ws.Range("F7").Value = 597
$ws.Range("F9").Value = 9053
$ws.Range("F10").Value = 1303
$ws.Range("F11").Value = 118
$ws.Range("F12").Value = 19

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 2067
$ws.Range("F4").Value = 7902
$ws.Range("F5").Value = 597
$ws.Range("F8").Value = 1303
$ws.Range("F9").Value = 118
$ws.Range("F10").Value = 4
$ws.Range("F12").Value = 160
$ws.Range("F13").Value = 1060
$ws.Range("F14").Value = 8536
$ws.Range("F18").Value = 1043
$ws.Range("F19").Value = 1018
$ws.Range("F20").Value = 1174
$ws.Range("F21").Value = 589
$ws.Range("F22").Value = 28
$ws.Range("F24").Value = 6
$ws.Range("F25").Value = 200
$ws.Range("F26").Value = 2
$ws.Range("F27").Value = 110
$ws.Range("F29").Value = 371
$ws.Range("F33").Value = 3513
$ws.Range("F34").Value = 929
$ws.Range("F36").Value = 504
$ws.Range("F37").Value = 40
$ws.Range("F39").Value = 109
$ws.Range("F40").Value = 639
$ws.Range("F43").Value = 46
$ws.Range("F45").Value = 22
